$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Fixed date placeholder text: 7/10/2019 -> 7/22/2019
#    Present on the slide master and on every slide layout's
#    "Date Placeholder" shape.
# ---------------------------------------------------------------------------
function Update-DateShapes($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq "7/10/2019") {
                $shp.TextFrame.TextRange.Text = "7/22/2019"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    Update-DateShapes $master.CustomLayouts.Item($li)
}

# ---------------------------------------------------------------------------
# 2) Slide 2 ("Git workflow"), TextBox 2 -- two paragraph edits.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(2)
$box = $slide.Shapes.Item(2)
$tr = $box.TextFrame.TextRange

# 2a) "... fetch the 'develop' branch to the local repo" becomes three runs:
#     " ... local " / "repo (You need go to the " / "repository directory first)"
$full = $tr.Text
$idx = $full.IndexOf("local repo")
$wordStart = $idx + "local ".Length + 1   # 1-based start of the word "repo"
$sub = $tr.Characters($wordStart, "repo".Length)
$sub.Text = "repo (You need go to the repository directory first)"

$full2 = $box.TextFrame.TextRange.Text
$tailText = "repository directory first)"
$tailPos = $full2.IndexOf($tailText) + 1   # 1-based
$subTail = $box.TextFrame.TextRange.Characters($tailPos, $tailText.Length)
$subTail.Text = $tailText

# 2b) "git" + " " + "pull " + "origin dev       " + "# (pull the recent from
#     remote repo update first)" collapses into "git" + a single merged run.
$full3 = $box.TextFrame.TextRange.Text
$idxGit = $full3.IndexOf("git pull origin dev")
$mergeStart = $idxGit + "git".Length + 1   # 1-based, right after "git"
$endMarker = "update first)"
$idxEnd = $full3.IndexOf($endMarker) + $endMarker.Length
$mergeLength = $idxEnd - ($idxGit + "git".Length)
$mergedText = $full3.Substring($idxGit + "git".Length, $mergeLength)
$mergeRange = $box.TextFrame.TextRange.Characters($mergeStart, $mergeLength)
$mergeRange.Text = $mergedText
